# "version final sin errores"
#
# Changes applied:
#  1. Rename the "Include from ICD-10" worksheet to "Include from Códigos de ejemp".
#  2. On the "Metadata" sheet:
#       - bump the Version value (row with "Version" in col A) from 0.4.0 to 0.7.0
#       - remove the "Jurisdiction"/"Chile" row entirely (rows below shift up)
#  3. On the "Include from ..." sheet, update the System URI value (row 4, col B)
#     to the new CodeSystem URL.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsCodes = $wb.Worksheets.Item(2)

# 1. Rename the second worksheet.
$wsCodes.Name = "Include from Códigos de ejemp"

# 2a. Update the Version value on the Metadata sheet.
$wsMeta.Range("B3").Value = "0.7.0"

# 2b. Delete the "Jurisdiction" / "Chile" row (row 11) entirely, shifting
#     everything below it up by one row.
$wsMeta.Rows.Item(11).Delete()

# 3. Update the System URI value on the codes sheet (row 4, column B).
$wsCodes.Range("B4").Value = "https://hospitallaflorida.cl/fhir/hlfhosp/CodeSystem/CSCie10Prueba"
